$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last updated" timestamp note in A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 03:50"

# 2. Update Estados Unidos (row 4) stats
$ws.Range("B4").Value = 188524
$ws.Range("C4").Value = 24736
$ws.Range("D4").Value = 7251
$ws.Range("E4").Value = 177384
$ws.Range("F4").Value = 4576
$ws.Range("G4").Value = 748
$ws.Range("H4").Value = 3889

# 3. Insert "Japon" as a new entry right after "Rumania" (row 34), shifting
#    Luxemburgo and Filipinas down one row (Pakistan and below stay put).
$ws.Range("A34").Value = "Japon"
$ws.Range("B34").Value = 2229
$ws.Range("C34").Value = 276
$ws.Range("D34").Value = 424
$ws.Range("E34").Value = 1739
$ws.Range("F34").Value = 69
$ws.Range("G34").Value = 10
$ws.Range("H34").Value = 66

$ws.Range("A35").Value = "Luxemburgo"
$ws.Range("B35").Value = 2178
$ws.Range("C35").Value = 190
$ws.Range("D35").Value = 80
$ws.Range("E35").Value = 2075
$ws.Range("F35").Value = 31
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 23

$ws.Range("A36").Value = "Filipinas"
$ws.Range("B36").Value = 2084
$ws.Range("C36").Value = 538
$ws.Range("D36").Value = 49
$ws.Range("E36").Value = 1947
$ws.Range("F36").Value = 1
$ws.Range("G36").Value = 10
$ws.Range("H36").Value = 88

# 4. Update Hong Kong (row 59) stats
$ws.Range("B59").Value = 715
$ws.Range("C59").Value = 32
$ws.Range("E59").Value = 583
